{"js": "// Replace the date heading and every two-digit x two-digit multiplication\n// \"problem=answer\" cell in the practice-sheet table with the next day's\n// values, per the commit's regenerated answer key.\nconst replacements = [\n  [\"2026-02-15 Sunday\", \"2026-02-16 Monday\"],\n  [\"23\u00d742=966\", \"66\u00d737=2442\"],\n  [\"44\u00d751=2244\", \"20\u00d799=1980\"],\n  [\"98\u00d755=5390\", \"80\u00d742=3360\"],\n  [\"45\u00d750=2250\", \"18\u00d749=882\"],\n  [\"47\u00d778=3666\", \"56\u00d739=2184\"],\n  [\"17\u00d760=1020\", \"98\u00d772=7056\"],\n  [\"53\u00d735=1855\", \"30\u00d775=2250\"],\n  [\"33\u00d777=2541\", \"49\u00d789=4361\"],\n  [\"34\u00d741=1394\", \"83\u00d759=4897\"],\n  [\"98\u00d727=2646\", \"74\u00d754=3996\"],\n  [\"98\u00d726=2548\", \"35\u00d792=3220\"],\n  [\"64\u00d729=1856\", \"80\u00d744=3520\"],\n  [\"87\u00d752=4524\", \"96\u00d758=5568\"],\n  [\"59\u00d758=3422\", \"48\u00d759=2832\"],\n  [\"87\u00d781=7047\", \"55\u00d757=3135\"],\n  [\"83\u00d756=4648\", \"45\u00d769=3105\"],\n  [\"20\u00d782=1640\", \"38\u00d758=2204\"],\n  [\"64\u00d755=3520\", \"66\u00d745=2970\"],\n  [\"99\u00d730=2970\", \"15\u00d753=795\"],\n  [\"60\u00d743=2580\", \"37\u00d786=3182\"],\n  [\"33\u00d747=1551\", \"95\u00d799=9405\"],\n  [\"25\u00d793=2325\", \"61\u00d762=3782\"],\n  [\"47\u00d734=1598\", \"97\u00d722=2134\"],\n  [\"19\u00d775=1425\", \"11\u00d785=935\"],\n  [\"87\u00d751=4437\", \"15\u00d776=1140\"],\n];\n\nconst body = context.document.body;\n\nfor (const [from, to] of replacements) {\n  // Each \"from\" string is unique in the document, so a literal (non-wildcard)\n  // search reliably returns exactly the single run that needs updating.\n  const results = body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for: ${from}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(to, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date heading and every two-digit x two-digit multiplication\n# \"problem=answer\" cell in the practice-sheet table with the next day's\n# values, per the commit's regenerated answer key.\n$d = $word.ActiveDocument\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$pairs = @(\n    @(\"2026-02-15 Sunday\", \"2026-02-16 Monday\"),\n    @(\"23\u00d742=966\", \"66\u00d737=2442\"),\n    @(\"44\u00d751=2244\", \"20\u00d799=1980\"),\n    @(\"98\u00d755=5390\", \"80\u00d742=3360\"),\n    @(\"45\u00d750=2250\", \"18\u00d749=882\"),\n    @(\"47\u00d778=3666\", \"56\u00d739=2184\"),\n    @(\"17\u00d760=1020\", \"98\u00d772=7056\"),\n    @(\"53\u00d735=1855\", \"30\u00d775=2250\"),\n    @(\"33\u00d777=2541\", \"49\u00d789=4361\"),\n    @(\"34\u00d741=1394\", \"83\u00d759=4897\"),\n    @(\"98\u00d727=2646\", \"74\u00d754=3996\"),\n    @(\"98\u00d726=2548\", \"35\u00d792=3220\"),\n    @(\"64\u00d729=1856\", \"80\u00d744=3520\"),\n    @(\"87\u00d752=4524\", \"96\u00d758=5568\"),\n    @(\"59\u00d758=3422\", \"48\u00d759=2832\"),\n    @(\"87\u00d781=7047\", \"55\u00d757=3135\"),\n    @(\"83\u00d756=4648\", \"45\u00d769=3105\"),\n    @(\"20\u00d782=1640\", \"38\u00d758=2204\"),\n    @(\"64\u00d755=3520\", \"66\u00d745=2970\"),\n    @(\"99\u00d730=2970\", \"15\u00d753=795\"),\n    @(\"60\u00d743=2580\", \"37\u00d786=3182\"),\n    @(\"33\u00d747=1551\", \"95\u00d799=9405\"),\n    @(\"25\u00d793=2325\", \"61\u00d762=3782\"),\n    @(\"47\u00d734=1598\", \"97\u00d722=2134\"),\n    @(\"19\u00d775=1425\", \"11\u00d785=935\"),\n    @(\"87\u00d751=4437\", \"15\u00d776=1140\"),\n)\n\nforeach ($pair in $pairs) {\n    $searchText = $pair[0]\n    $replaceText = $pair[1]\n\n    # Each $searchText is unique in the document, so searching the whole\n    # document body and replacing all (one) match is safe and precise.\n    $find = $d.Content.Find\n    $found = $find.Execute($searchText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll)\n    if (-not $found) {\n        throw \"No match found for: $searchText\"\n    }\n}\n"}
